$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Calc")

# A new Morgan federal poll has come in, becoming "Latest Morgan" (row 6).
# The previous "Latest Morgan" (row 6) becomes "Second Morgan" (row 7), and
# the previous "Second Morgan" (row 7) becomes "Third Morgan" (row 8); the
# former "Third Morgan" (old row 8) data is dropped. Update bottom rows
# first so the literal values below (captured from the prior state) are
# written correctly regardless of execution order.

# Row 8 (Third Morgan) <- values previously in row 7 (Second Morgan)
$ws.Range("B8").Value = 56.5
$ws.Range("C8").Value = 56.5
$ws.Range("D8").Value = 60
$ws.Range("E8").Value = 50
$ws.Range("F8").Value = 52
$ws.Range("G8").Value = 59.5

# Row 7 (Second Morgan) <- values previously in row 6 (Latest Morgan)
$ws.Range("B7").Value = 56.5
$ws.Range("C7").Value = 56
$ws.Range("D7").Value = 63.5
$ws.Range("E7").Value = 48.5
$ws.Range("F7").Value = 52
$ws.Range("G7").Value = 52.5

# Row 6 (Latest Morgan) <- brand new SA YouGov / Morgan federal poll numbers
$ws.Range("B6").Value = 56
$ws.Range("C6").Value = 56.5
$ws.Range("D6").Value = 60
$ws.Range("E6").Value = 48
$ws.Range("F6").Value = 53
$ws.Range("G6").Value = 54.5

# Update the view: scroll position reset and selection moved to G7
$ws.Activate()
$ws.Range("G7").Select()
